# Apply the commit's edit:
#  - Sheet1: no data change, just move the remembered selection from B4 to D2.
#  - Api (sheet2): append a new data row (row 4) -> API Data01 / 6046 / 1739,
#    growing the used range from A1:C3 to A1:C4, and leave the sheet's
#    selection at the top-left (A1) default.

$wb = $excel.ActiveWorkbook

# --- Sheet1: move the saved selection to D2 -------------------------------
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Activate() | Out-Null
$ws1.Range("D2").Select() | Out-Null

# --- Api: add the new "API Data01 / 6046 / 1739" row ----------------------
$ws2 = $wb.Worksheets.Item("Api")
$ws2.Activate() | Out-Null

$ws2.Range("A4").Value = "API Data01"
$ws2.Range("B4").Value = "6046"
$ws2.Range("C4").Value = "1739"

# Leave the Api sheet's selection back at the default top-left cell.
$ws2.Range("A1").Select() | Out-Null
